$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.177.21"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.271.41"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.34%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -3.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.06"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.48%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.63"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.621.47"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.259.58"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.773"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.150.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.96"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.68"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("E25").Value = "  -4.80%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.31"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.01"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("E36").Value = "  -5.64%  "
$ws.Range("E37").Value = "  -5.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.02"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -9.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0988"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.40%  "
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.18%  "
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.957.65"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.494.02"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.81"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.58%  "
